# Update "Datos actualizados" timestamp (A1)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

$ws.Range("A1").Value = "Datos actualizados a 16 de Julio de 2020 a las 05:33"

# Rows 130/131: Suazilandia overtakes Ruanda in the ranking.
$ws.Range("A130").Value = "Suazilandia"
$ws.Range("B130").Value = 1489
$ws.Range("D130").Value = 724
$ws.Range("E130").Value = 745
$ws.Range("H130").Value = 20

$ws.Range("A131").Value = "Ruanda"
$ws.Range("B131").Value = 1435
$ws.Range("D131").Value = 752
$ws.Range("E131").Value = 679
$ws.Range("H131").Value = 4

# Row 138: Niger small update
$ws.Range("B138").Value = 1100
$ws.Range("D138").Value = 993
$ws.Range("E138").Value = 38
$ws.Range("H138").Value = 69

# Rows 140-143: Liberia surges and overtakes Uganda, Burkina Faso and Chipre.
$ws.Range("A140").Value = "Liberia"
$ws.Range("B140").Value = 1056
$ws.Range("D140").Value = 447
$ws.Range("E140").Value = 558
$ws.Range("H140").Value = 51

$ws.Range("A141").Value = "Uganda"
$ws.Range("B141").Value = 1043
$ws.Range("D141").Value = 1004
$ws.Range("E141").Value = 39
$ws.Range("H141").Value = 0

$ws.Range("A142").Value = "Burkina Faso"
$ws.Range("B142").Value = 1038
$ws.Range("D142").Value = 882
$ws.Range("E142").Value = 103
$ws.Range("H142").Value = 53

$ws.Range("A143").Value = "Republica de Chipre"
$ws.Range("B143").Value = 1025
$ws.Range("D143").Value = 839
$ws.Range("E143").Value = 167
$ws.Range("H143").Value = 19

# Row 148: Principado de Andorra small update
$ws.Range("B148").Value = 862
$ws.Range("E148").Value = 7

# Row 152: Santo Tome y Principe update
$ws.Range("B152").Value = 737
$ws.Range("D152").Value = 322
$ws.Range("E152").Value = 401

# Rows 209/210: Groenlandia and Islas Malvinas swap (tied values, label-only change)
$ws.Range("A209").Value = "Groenlandia"
$ws.Range("A210").Value = "Islas Malvinas"

Write-Host "Applied country/province updates"
